$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks (stored XML width = ColumnWidth + 5/6 char padding) ---
# B: 8 -> 7
$ws.Columns("B").ColumnWidth = 6.166666666666667
# F: 8 -> 7
$ws.Columns("F").ColumnWidth = 6.166666666666667
# P: 8 -> 7
$ws.Columns("P").ColumnWidth = 6.166666666666667
# Q: 7 -> 8
$ws.Columns("Q").ColumnWidth = 7.166666666666667
# W: 8 -> 7
$ws.Columns("W").ColumnWidth = 6.166666666666667

# --- Refresh data rows 2-5 with the new sample window (custom accuracy run) ---
$row2 = New-Object "object[,]" 1,34
$row2[0,0] = 45161.50694444445
$row2[0,1] = 4.639
$row2[0,2] = 3.949
$row2[0,3] = 0
$row2[0,4] = 5.682
$row2[0,5] = 5.56
$row2[0,6] = 1.404
$row2[0,7] = 7.057
$row2[0,8] = 3.269
$row2[0,9] = 2.974
$row2[0,10] = 2.59
$row2[0,11] = 3.719
$row2[0,12] = 4.974
$row2[0,13] = 2.211
$row2[0,14] = 2.2
$row2[0,15] = 3.582
$row2[0,16] = 1.706
$row2[0,17] = 0.605
$row2[0,18] = 0.024
$row2[0,19] = 37.258
$row2[0,20] = 6.598
$row2[0,21] = 4.255
$row2[0,22] = 5.158
$row2[0,23] = 1.722
$row2[0,24] = 0.377
$row2[0,25] = 3.036
$row2[0,26] = 1.789
$row2[0,27] = 1.745
$row2[0,28] = 5.904
$row2[0,29] = 3.659
$row2[0,30] = 3.978
$row2[0,31] = 5.202
$row2[0,32] = 0.995
$row2[0,33] = 2.83
$ws.Range("A2:AH2").Value2 = $row2

$row3 = New-Object "object[,]" 1,34
$row3[0,0] = 45161.51388888889
$row3[0,1] = 18.14
$row3[0,2] = 13.763
$row3[0,3] = 0.434
$row3[0,4] = 37.666
$row3[0,5] = 31.64
$row3[0,6] = 13.257
$row3[0,7] = 47.795
$row3[0,8] = 21.029
$row3[0,9] = 10.26
$row3[0,10] = 14.331
$row3[0,11] = 15.589
$row3[0,12] = 16.915
$row3[0,13] = 5.077
$row3[0,14] = 13.528
$row3[0,15] = 19.76
$row3[0,16] = 11.034
$row3[0,17] = 0.594
$row3[0,18] = 0.428
$row3[0,19] = 203.985
$row3[0,20] = 38.126
$row3[0,21] = 13.332
$row3[0,22] = 26.16
$row3[0,23] = 13.287
$row3[0,24] = 1.821
$row3[0,25] = 24.429
$row3[0,26] = 11.146
$row3[0,27] = 9.978999999999999
$row3[0,28] = 12.962
$row3[0,29] = 16.393
$row3[0,30] = 1.327
$row3[0,31] = 43.158
$row3[0,32] = 7.064
$row3[0,33] = 15.728
$ws.Range("A3:AH3").Value2 = $row3

$row4 = New-Object "object[,]" 1,34
$row4[0,0] = 45161.52083333334
$row4[0,1] = 2.638
$row4[0,2] = 2.058
$row4[0,3] = 0.007
$row4[0,4] = 4.599
$row4[0,5] = 4.128
$row4[0,6] = 1.421
$row4[0,7] = 13.644
$row4[0,8] = 2.616
$row4[0,9] = 1.841
$row4[0,10] = 1.904
$row4[0,11] = 2.127
$row4[0,12] = 2.503
$row4[0,13] = 1.038
$row4[0,14] = 1.632
$row4[0,15] = 2.804
$row4[0,16] = 1.202
$row4[0,17] = 0.346
$row4[0,18] = 0.013
$row4[0,19] = 21.143
$row4[0,20] = 5.289
$row4[0,21] = 2.023
$row4[0,22] = 3.914
$row4[0,23] = 1.641
$row4[0,24] = 0.245
$row4[0,25] = 6.354
$row4[0,26] = 1.437
$row4[0,27] = 1.405
$row4[0,28] = 2.334
$row4[0,29] = 2.243
$row4[0,30] = 0.796
$row4[0,31] = 12.998
$row4[0,32] = 0.8120000000000001
$row4[0,33] = 1.958
$ws.Range("A4:AH4").Value2 = $row4

$row5 = New-Object "object[,]" 1,34
$row5[0,0] = 45161.52777777778
$row5[0,1] = 11.22
$row5[0,2] = 8.48
$row5[0,3] = 0.29
$row5[0,4] = 23.6
$row5[0,5] = 19.74
$row5[0,6] = 8.35
$row5[0,7] = 31.13
$row5[0,8] = 13.17
$row5[0,9] = 6.31
$row5[0,10] = 8.949999999999999
$row5[0,11] = 9.66
$row5[0,12] = 10.41
$row5[0,13] = 3.09
$row5[0,14] = 8.470000000000001
$row5[0,15] = 12.31
$row5[0,16] = 6.88
$row5[0,17] = 0.36
$row5[0,18] = 0.27
$row5[0,19] = 124.07
$row5[0,20] = 23.71
$row5[0,21] = 8.18
$row5[0,22] = 16.18
$row5[0,23] = 8.34
$row5[0,24] = 1.13
$row5[0,25] = 15.34
$row5[0,26] = 6.99
$row5[0,27] = 6.25
$row5[0,28] = 7.82
$row5[0,29] = 10.16
$row5[0,30] = 0.57
$row5[0,31] = 27.89
$row5[0,32] = 4.42
$row5[0,33] = 9.81
$ws.Range("A5:AH5").Value2 = $row5

# --- Dataset now spans only 4 data rows; drop the stale 6th row ---
$ws.Rows(6).Delete()
